$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210, shifting existing rows 210-278 down to 211-279
$ws.Rows.Item(210).Insert()

# Populate the new row 210 with the new data record
$ws.Cells.Item(210, 1).Value = 4
$ws.Cells.Item(210, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(210, 3).Value = "Los Lagos"
$ws.Cells.Item(210, 4).Value = 44985
$ws.Cells.Item(210, 5).Value = 10
$ws.Cells.Item(210, 6).Value = 100112009
$ws.Cells.Item(210, 7).Value = "Acelga"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 80
$ws.Cells.Item(210, 11).Value = 10000
$ws.Cells.Item(210, 12).Value = 10000
$ws.Cells.Item(210, 13).Value = 10000
$ws.Cells.Item(210, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(210, 15).Value = "Región Metropolitana"
$ws.Cells.Item(210, 16).Value = 833
$ws.Cells.Item(210, 17).Value = 12
$ws.Cells.Item(210, 18).Value = "Hortaliza"
